$d = $word.ActiveDocument

$replacements = @(
    @{old="2025-01-29 Wednesday"; new="2025-01-30 Thursday"},
    @{old="435÷4=108, 3"; new="832÷2=416, 0"},
    @{old="384÷6=64, 0"; new="890÷2=445, 0"},
    @{old="764÷2=382, 0"; new="845÷7=120, 5"},
    @{old="282÷2=141, 0"; new="870÷6=145, 0"},
    @{old="318÷4=79, 2"; new="807÷8=100, 7"},
    @{old="658÷5=131, 3"; new="113÷7=16, 1"},
    @{old="491÷6=81, 5"; new="366÷4=91, 2"},
    @{old="237÷7=33, 6"; new="556÷3=185, 1"},
    @{old="589÷6=98, 1"; new="708÷6=118, 0"},
    @{old="130÷5=26, 0"; new="319÷5=63, 4"},
    @{old="811÷7=115, 6"; new="164÷2=82, 0"},
    @{old="249÷5=49, 4"; new="504÷9=56, 0"},
    @{old="406÷3=135, 1"; new="942÷3=314, 0"},
    @{old="530÷7=75, 5"; new="349÷5=69, 4"},
    @{old="732÷3=244, 0"; new="131÷4=32, 3"},
    @{old="819÷6=136, 3"; new="611÷9=67, 8"},
    @{old="574÷8=71, 6"; new="437÷4=109, 1"},
    @{old="316÷3=105, 1"; new="982÷9=109, 1"},
    @{old="968÷3=322, 2"; new="320÷2=160, 0"},
    @{old="498÷6=83, 0"; new="400÷7=57, 1"},
    @{old="462÷8=57, 6"; new="779÷7=111, 2"},
    @{old="933÷9=103, 6"; new="396÷5=79, 1"},
    @{old="654÷4=163, 2"; new="991÷2=495, 1"},
    @{old="349÷7=49, 6"; new="938÷6=156, 2"},
    @{old="604÷5=120, 4"; new="705÷8=88, 1"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $r.new, 2)
}

Write-Output "Done"
